# "refactor, meet initial requirements"
#
# Adds a computed "percentage" column (hours as a share of the PI's /
# department's total hours, expressed as a percentage) to both sheets:
#   - "PI hours"   (sheet 1): name, hours, percentage, dept
#   - "dept hours" (sheet 2): dept, hours, percentage
#
# On sheet 1 the new column is inserted between "hours" and "dept" (so
# "dept" shifts from column D to column E); on sheet 2 it is appended
# after "hours".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "PI hours" -> insert "percentage" before "dept"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Shift the existing "dept" column (D) one place to the right; this
# carries the header's border/bold/centered style along with it, so the
# new column D inherits that same formatting automatically.
$ws1.Columns.Item(4).Insert()

$ws1.Cells.Item(1, 4).Value = "percentage"
$ws1.Cells.Item(2, 4).Value = 87.5
$ws1.Cells.Item(3, 4).Value = 12.5

# ---------------------------------------------------------------------
# Sheet 2: "dept hours" -> append "percentage" after "hours"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Copy the "hours" header's formatting (bold, centered, thin border)
# onto the new header cell before writing its text.
$ws2.Cells.Item(1, 2).Copy()
$ws2.Cells.Item(1, 4).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Cells.Item(1, 4).Value = "percentage"
$ws2.Cells.Item(2, 4).Value = 47.05882352941177
$ws2.Cells.Item(3, 4).Value = 41.1764705882353
$ws2.Cells.Item(4, 4).Value = 5.882352941176471
$ws2.Cells.Item(5, 4).Value = 5.882352941176471
